# Updates the cryptos list (Price / Volume(1h) columns) for the Fri Jun 21
# 22:09:13 UTC 2024 GitHub Actions refresh.
#
# Note: several "Price" values look like plain decimals (e.g. "1.00",
# "27.54") which Excel would otherwise auto-convert to numbers on
# assignment, silently dropping the trailing zero / text formatting that
# the source sheet relies on (it stores every Price/Volume cell as text).
# Prefixing those values with a leading apostrophe forces Excel to keep
# them as literal text, matching the original workbook's cell typing.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.114.62'
$ws.Range('E2').Value = '  -1.47%  '
$ws.Range('D3').Value = '3.517.76'
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''586.80'
$ws.Range('E5').Value = '  -0.97%  '
$ws.Range('D6').Value = '''133.53'
$ws.Range('E6').Value = '  -0.33%  '
$ws.Range('D7').Value = '3.518.81'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D9').Value = '''0.489'
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('D11').Value = '''7.14'
$ws.Range('E11').Value = '  -0.22%  '
$ws.Range('E12').Value = '  -2.36%  '
$ws.Range('D13').Value = '4.115.00'
$ws.Range('E13').Value = '  +0.01%  '
$ws.Range('D14').Value = '''27.54'
$ws.Range('E14').Value = '  -0.61%  '
$ws.Range('E15').Value = '  +1.40%  '
$ws.Range('D16').Value = '''0.0000179'
$ws.Range('E16').Value = '  -1.59%  '
$ws.Range('D17').Value = '3.498.32'
$ws.Range('E17').Value = '  -0.57%  '
$ws.Range('D18').Value = '64.143.10'
$ws.Range('E18').Value = '  -1.37%  '
$ws.Range('D19').Value = '''9.83'
$ws.Range('E19').Value = '  -2.40%  '
$ws.Range('D20').Value = '''13.88'
$ws.Range('E20').Value = '  -3.12%  '
$ws.Range('E21').Value = '  -1.01%  '
$ws.Range('D22').Value = '''382.81'
$ws.Range('E22').Value = '  -2.32%  '
$ws.Range('E23').Value = '  -1.25%  '
$ws.Range('D24').Value = '3.659.32'
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('D25').Value = '''73.93'
$ws.Range('E25').Value = '  -0.86%  '
$ws.Range('D26').Value = '''1.00'
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('D27').Value = '''5.62'
$ws.Range('E27').Value = '  -1.29%  '
$ws.Range('E28').Value = '  +2.78%  '
$ws.Range('D29').Value = '''1.58'
$ws.Range('E29').Value = '  -1.58%  '
$ws.Range('D30').Value = '''7.47'
$ws.Range('E30').Value = '  -2.58%  '
$ws.Range('E31').Value = '  +0.09%  '
$ws.Range('D32').Value = '''8.47'
$ws.Range('E32').Value = '  +1.96%  '
$ws.Range('E33').Value = '  -1.64%  '
$ws.Range('D34').Value = '3.531.91'
$ws.Range('E34').Value = '  +0.19%  '
$ws.Range('E36').Value = '  -2.23%  '
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('D38').Value = '''5.39'
$ws.Range('E38').Value = '  +2.45%  '
$ws.Range('E39').Value = '  -0.21%  '
$ws.Range('D40').Value = '''6.95'
$ws.Range('E40').Value = '  -0.21%  '
$ws.Range('D41').Value = '''160.23'
$ws.Range('E41').Value = '  -4.64%  '
$ws.Range('D42').Value = '''0.0786'
$ws.Range('E42').Value = '  -2.52%  '
$ws.Range('D43').Value = '''26.70'
$ws.Range('E43').Value = '  +3.30%  '
$ws.Range('D44').Value = '''0.813'
$ws.Range('E44').Value = '  -0.82%  '
$ws.Range('D45').Value = '''1.00'
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('D46').Value = '''1.22'
$ws.Range('E46').Value = '  -3.15%  '
$ws.Range('D47').Value = '''41.64'
$ws.Range('E47').Value = '  -3.08%  '
$ws.Range('D48').Value = '''4.41'
$ws.Range('E48').Value = '  -0.46%  '
$ws.Range('E49').Value = '  -3.04%  '
$ws.Range('D50').Value = '2.483.85'
$ws.Range('E50').Value = '  +2.48%  '
$ws.Range('E51').Value = '  -1.35%  '
